# Insert a new "2022-Q3" sheet with fund-holding detail data, and update the
# "总计" (summary) sheet to include the new quarter and re-sequence existing rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value to a cell while forcing "Text" number format so that
# numeric-looking strings (fund codes with leading zeros, figures with
# trailing zeros, etc.) are preserved exactly as text instead of being
# coerced into numbers.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet named "2022-Q3" immediately before "2022-Q2".
#    All the existing quarterly sheets keep their own names and simply shift
#    one position to the right.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($refSheet)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with the fund holdings table.
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $cell = $q3.Cells.Item(1, $c + 2)
    $cell.Value = $headers[$c]
}
# Reuse the existing bold/bordered header style (style index used by B1 on
# the "2022-Q2" sheet) instead of creating a brand-new style.
$refSheet.Cells.Item(1, 2).Copy() | Out-Null
$q3.Range($q3.Cells.Item(1, 2), $q3.Cells.Item(1, 8)).PasteSpecial(-4122) | Out-Null

$fundRows = @(
    @(0, '000021', '华夏优势增长混合', '55.95', '89.84', '1.84', '1.0295', 9),
    @(1, '010180', '华夏科技龙头两年定期开放混合', '18.82', '93.38', '4.68', '0.8808', 3),
    @(2, '000061', '华夏盛世混合', '14.10', '82.60', '2.86', '0.4033', 4),
    @(3, '001042', '华夏领先股票', '11.44', '89.89', '2.69', '0.3077', 8),
    @(4, '012173', '国泰兴泽优选一年持有期混合A', '8.41', '88.23', '3.12', '0.2624', 9),
    @(5, '012174', '国泰兴泽优选一年持有期混合C', '6.17', '88.23', '3.12', '0.1925', 9),
    @(6, '001924', '华夏国企改革灵活配置混合', '2.66', '88.64', '4.82', '0.1282', 9),
    @(7, '010016', '华夏科技前沿6个月定期开放混合A', '4.76', '84.35', '2.43', '0.1157', 4),
    @(8, '010017', '华夏科技前沿6个月定期开放混合C', '2.05', '84.35', '2.43', '0.0498', 4),
    @(9, '005083', '诺德量化蓝筹增强混合C', '0.57', '92.85', '2.55', '0.0145', 9),
    @(10, '002292', '诺安益鑫灵活配置混合A', '0.39', '61.16', '3.32', '0.0129', 9),
    @(11, '003238', '新华外延增长主题灵活配置混合', '0.50', '57.43', '1.66', '0.0083', 10),
    @(12, '015466', '太平中证1000指数增强A', '0.37', '92.23', '1.11', '0.0041', 2),
    @(13, '014550', '诺安益鑫灵活配置混合C', '0.02', '61.16', '3.32', '0.0007', 9),
    @(14, '015467', '太平中证1000指数增强C', '0.02', '92.23', '1.11', '0.0002', 2),
    @(15, '005082', '诺德量化蓝筹增强混合A', '0.00', '92.85', '2.55', $null, 9)
)

foreach ($row in $fundRows) {
    $r = [int]$row[0] + 2
    $q3.Cells.Item($r, 1).Value = [int]$row[0]
    Set-TextValue $q3.Cells.Item($r, 2) $row[1]
    Set-TextValue $q3.Cells.Item($r, 3) $row[2]
    Set-TextValue $q3.Cells.Item($r, 4) $row[3]
    Set-TextValue $q3.Cells.Item($r, 5) $row[4]
    Set-TextValue $q3.Cells.Item($r, 6) $row[5]
    if ($row[6] -eq $null) {
        $q3.Cells.Item($r, 7).Value = 0
    } else {
        Set-TextValue $q3.Cells.Item($r, 7) $row[6]
    }
    $q3.Cells.Item($r, 8).Value = [int]$row[7]
}
# Reuse the same style used for column A on the reference sheet for the new
# index column.
$refSheet.Cells.Item(2, 1).Copy() | Out-Null
$q3.Range($q3.Cells.Item(2, 1), $q3.Cells.Item(17, 1)).PasteSpecial(-4122) | Out-Null
foreach ($row in $fundRows) {
    $r = [int]$row[0] + 2
    $q3.Cells.Item($r, 1).Value = [int]$row[0]
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: add the 2022-Q3 row at the top of the
#    data and push the older quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @("2022-Q3", 16, 3.41),
    @("2022-Q2", 11, 1.78),
    @("2022-Q1", 17, 11.46),
    @("2021-Q4", 28, 13.66),
    @("2021-Q3", 12, 5.75),
    @("2021-Q2", 6, 2.94)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $summaryRows[$i][0]
    $total.Cells.Item($r, 3).Value = $summaryRows[$i][1]
    $total.Cells.Item($r, 4).Value = $summaryRows[$i][2]
}
# Make sure the newly created row (row 7, for 2021-Q2) carries the same
# style as the rest of the index column.
$total.Cells.Item(6, 1).Copy() | Out-Null
$total.Cells.Item(7, 1).PasteSpecial(-4122) | Out-Null
$total.Cells.Item(7, 1).Value = 5
